$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "健康増進成分: Mystic Spice Chai Tea の各成分は、その自然な健康上の利点のために選択されています。" "健康に良い素材: Mystic Spice Chai Tea には、健康効果を考慮して選ばれた天然素材が使用されています。"

Replace-Text "豊かな香りと味: 私たちのチャイは暖かく、辛い香りと深く、活気のある味わいを持ち、一日を始めたり、夜にリラックスするのに最適な飲み物です。" "豊かな香りと風味: 温かくスパイシーな香りと深く爽快な味わいで、一日の始まりや夜のリラックスタイムに最適な飲み物です。"

Replace-Text "どんな淹れ方でも: ホット ティー、さわやかなアイス ティー、クリーミーなラテなど、私たちのブレンドはどんな好みにも合うように作られています。" "多様な楽しみ方: ホットだけでなく、爽やかなアイスティーやクリーミーなラテなど、好みに合わせて自由にお楽しみいただける万能なブレンドです。"

Replace-Text "サステナブル ソース: 持続可能性に努め、有機農業を実践する小規模農場から原料を調達し、最高の品質だけでなく、地球の福祉も確保しています。" "持続可能な調達: 持続可能性を重視し、有機農業を実践する小規模農場から原料を調達することで、最高品質を実現するだけでなく地球環境にも配慮しています。"

Replace-Text "エレガントな包装: Mystic Spice Chai Tea は美しくデザインされ、環境に優しい包装で提供され、紅茶愛好家にふさわしい贈り物や自分への豪華なご褒美になります。" "エレガントなパッケージ: Mystic Spice Chai Tea は、環境に配慮した美しいデザインのパッケージに梱包されているため、お茶が好きな方に贈るギフトや自分自身への贅沢なご褒美として最適です。"

Replace-Text "顧客満足度保証: 私たちは製品を支持し、満足度を保証します。" "顧客満足度保証: 当社は製品の品質に自信を持っており、満足度保証を提供しています。"

Replace-Text "中度" "中"

Replace-Text "ただし、プロモーション計画と戦略は、市場の状況の変化や顧客からのフィードバックに応じて常に監視、評価、調整する必要があります。" "ただし、販売推進計画と戦略は、市場の状況の変化や顧客からのフィードバックに応じて常に監視、評価、調整する必要があります。"
